{"js": "// Highlight the word \"Populations\" in yellow within the\n// \"\u2022 Populations \u2013 EA optimiser\" bullet paragraph.\nconst body = context.document.body;\nconst results = body.search(\"Populations\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].font.highlightColor = \"Yellow\";\n}\n\nawait context.sync();\n", "ps1": "# Highlight the word \"Populations\" in yellow within the\n# \"\u2022 Populations \u2013 EA optimiser\" bullet paragraph.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Populations\"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $true\n\nif ($rng.Find.Execute()) {\n    # wdYellow = 7\n    $rng.Font.HighlightColorIndex = 7\n}\n"}
